# Add a new program row (TUM_MSNE) to the Program_choosing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

$ws.Range("A7").Value = "TUM_MSNE"
$ws.Range("B7").Value = "Yes"

# Extend the Yes/No list validation from B1:B6 down to B1:B7.
$ws.Range("B1:B7").Validation.Delete()
$ws.Range("B1:B7").Validation.Add(3, 1, 1, """Yes,No""")
$ws.Range("B1:B7").Validation.IgnoreBlank = $true
$ws.Range("B1:B7").Validation.InCellDropdown = $true
$ws.Range("B1:B7").Validation.ShowInput = $true
$ws.Range("B1:B7").Validation.ShowError = $true

# Move the active selection to the newly-added cell B7.
$ws.Range("B7").Select() | Out-Null
